$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Fri Jun 16 19:28:49 UTC 2023 with GitHub Actions

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.311.12"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.716.49"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.58"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4697"
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2631"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06218"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.710.75"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07068"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5883"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.24"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.291.80"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006812"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.929.31"
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.556"
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.799"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.335"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.28"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.15"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.406"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.759"
$ws.Range("E28").Value = "  +4.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.80"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.036"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.686"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07737"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04415"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.613"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6196"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9683"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9197"
$ws.Range("E37").Value = "  +8.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "112.67"
$ws.Range("E38").Value = "  +14.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.407"
$ws.Range("E39").Value = "  -7.98%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.895"
$ws.Range("E41").Value = "  +4.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01462"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.255"
$ws.Range("E43").Value = "  +12.08%  "
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("E45").Value = "  +3.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.232"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05290"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.635"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.218"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3370"
$ws.Range("E51").Value = "  +1.16%  "
